$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptocurrency price and volume(1h) figures
$ws.Range("D2").Value = '35.035.02'
$ws.Range("E2").Value = '  +1.24%  '
$ws.Range("D3").Value = '1.816.24'
$ws.Range("E3").Value = '  -0.21%  '
$ws.Range("E4").Value = '  +0.41%  '
$ws.Range("D5").Value = "'" + '233.31'
$ws.Range("E5").Value = '  +3.15%  '
$ws.Range("D6").Value = "'" + '0.617'
$ws.Range("E6").Value = '  +1.48%  '
$ws.Range("E7").Value = '  +0.49%  '
$ws.Range("D8").Value = "'" + '40.09'
$ws.Range("E8").Value = '  -10.64%  '
$ws.Range("D9").Value = "'" + '0.326'
$ws.Range("E9").Value = '  +9.54%  '
$ws.Range("D10").Value = "'" + '0.0685'
$ws.Range("E10").Value = '  +1.09%  '
$ws.Range("E11").Value = '  -0.27%  '
$ws.Range("E12").Value = '  -0.18%  '
$ws.Range("D13").Value = '1.836.28'
$ws.Range("E13").Value = '  +0.82%  '
$ws.Range("D14").Value = "'" + '11.11'
$ws.Range("E14").Value = '  -0.33%  '
$ws.Range("D15").Value = "'" + '4.69'
$ws.Range("E15").Value = '  +2.94%  '
$ws.Range("D16").Value = "'" + '0.661'
$ws.Range("E16").Value = '  +2.84%  '
$ws.Range("D17").Value = '34.974.70'
$ws.Range("E17").Value = '  +1.25%  '
$ws.Range("D18").Value = "'" + '69.56'
$ws.Range("E18").Value = '  +2.14%  '
$ws.Range("E19").Value = '  +1.16%  '
$ws.Range("D20").Value = "'" + '238.57'
$ws.Range("E20").Value = '  -1.40%  '
$ws.Range("D21").Value = "'" + '11.84'
$ws.Range("E21").Value = '  +1.34%  '
$ws.Range("D22").Value = "'" + '4.67'
$ws.Range("E22").Value = '  +1.96%  '
$ws.Range("E23").Value = '  +0.40%  '
$ws.Range("E24").Value = '  +4.14%  '
$ws.Range("D25").Value = "'" + '172.82'
$ws.Range("E25").Value = '  +1.12%  '
$ws.Range("D26").Value = "'" + '7.83'
$ws.Range("E26").Value = '  +0.01%  '
$ws.Range("D27").Value = "'" + '17.46'
$ws.Range("E27").Value = '  -1.62%  '
$ws.Range("E28").Value = '  -0.81%  '
$ws.Range("E29").Value = '  +31.82%  '
$ws.Range("E30").Value = '  +0.49%  '
$ws.Range("D31").Value = '3.338.66'
$ws.Range("E31").Value = '  +37.41%  '
$ws.Range("D32").Value = "'" + '0.0555'
$ws.Range("E32").Value = '  +6.56%  '
$ws.Range("D33").Value = "'" + '3.94'
$ws.Range("E33").Value = '  +1.90%  '
$ws.Range("D34").Value = "'" + '3.97'
$ws.Range("E34").Value = '  +0.65%  '
$ws.Range("E35").Value = '  -2.93%  '
$ws.Range("D36").Value = "'" + '93.22'
$ws.Range("E36").Value = '  +4.30%  '
$ws.Range("E37").Value = '  +6.71%  '
$ws.Range("D38").Value = "'" + '0.682'
$ws.Range("E38").Value = '  +3.34%  '
$ws.Range("E39").Value = '  +1.12%  '
$ws.Range("D40").Value = "'" + '1.29'
$ws.Range("E40").Value = '  +5.45%  '
$ws.Range("D41").Value = '1.307.77'
$ws.Range("E41").Value = '  -1.01%  '
$ws.Range("D42").Value = "'" + '0.985'
$ws.Range("E42").Value = '  +2.63%  '
$ws.Range("D43").Value = "'" + '2.34'
$ws.Range("E43").Value = '  -2.93%  '
$ws.Range("D44").Value = "'" + '14.62'
$ws.Range("E44").Value = '  -4.56%  '
$ws.Range("D45").Value = "'" + '2.46'
$ws.Range("E45").Value = '  +1.40%  '
$ws.Range("E46").Value = '  -1.43%  '
$ws.Range("D47").Value = "'" + '6.31'
$ws.Range("E47").Value = '  +6.92%  '
$ws.Range("D48").Value = "'" + '0.0512'
$ws.Range("E48").Value = '  -1.11%  '
$ws.Range("E49").Value = '  +0.64%  '
$ws.Range("D50").Value = "'" + '1.01'
$ws.Range("E50").Value = '  +0.43%  '
$ws.Range("E51").Value = '  +5.54%  '
